$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header column I: "COD ESPECIE BEM" ---
$ws.Range("I1").Value = "COD ESPECIE BEM"
$ws.Columns("I").ColumnWidth = 18.140625

# --- Row 2 edits: date + sequencia rateio change ---
$ws.Range("C2").Value = 45894
$ws.Range("E2").Value = 2
$ws.Range("I2").Value = 412

# --- New row 3 (mirrors row 2, with a few different values) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "9999999999-000"
$ws.Range("C3").Value = 45894
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 65
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = 1500.5
$ws.Range("I3").Value = 412

# --- Apply the Consolas (themed) font to the plain data cells ---
$ws.Range("A2:B2").Font.Name = "Consolas"
$ws.Range("D2:F2").Font.Name = "Consolas"
$ws.Range("I2").Font.Name = "Consolas"
$ws.Range("A3:B3").Font.Name = "Consolas"
$ws.Range("D3:F3").Font.Name = "Consolas"
$ws.Range("I3").Font.Name = "Consolas"

# --- Date + currency cells get the Consolas font too, plus their number formats ---
$ws.Range("C2").Font.Name = "Consolas"
$ws.Range("C2").NumberFormat = "m/d/yyyy"
$ws.Range("C3").Font.Name = "Consolas"
$ws.Range("C3").NumberFormat = "m/d/yyyy"

$ws.Range("H2").Font.Name = "Consolas"
$ws.Range("H2").NumberFormat = "#,##0.00"
$ws.Range("H3").Font.Name = "Consolas"
$ws.Range("H3").NumberFormat = "#,##0.00"

# G2 keeps its existing #,##0.0000 Consolas look; G3 gets a 4-decimal (no thousands) look
$ws.Range("G3").Font.Name = "Consolas"
$ws.Range("G3").NumberFormat = "0.0000"

# --- Extra formatted-but-empty rows below (4-7) ---
$ws.Range("A4:H7").Font.Name = "Consolas"

# --- Stray formatting left on a couple of far-flung empty cells ---
$ws.Range("D9").Font.Name = "Calibri"
$ws.Range("G10").Font.Name = "Calibri"
$ws.Range("J5").Font.Underline = $true

# --- View/selection bookkeeping ---
$ws.Range("E3").Select
